$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 03.02.2022 11:15"

# Row 3 (Tesco) updates:
# B3 gets a new current price, C3 takes the previous B3 value (old price).
$ws.Range("B3").Value = 35.5
$ws.Range("C3").Value = 34.5

# D3 becomes the text delta "+1.0" (must stay literal text, not be
# auto-coerced into the number 1). Format the cell as Text first so the
# leading "+" sign is preserved, then reset the cell style back to Normal
# (the target keeps no explicit style index on this cell).
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+1.0"
$ws.Range("D3").Style = "Normal"

# E3 becomes a plain text timestamp string (it previously held a numeric
# Excel date serial with a date-format style applied). Assign the text and
# drop the old date style back to Normal/default.
$ws.Range("E3").Value = "2022-02-03 11:17:15"
$ws.Range("E3").Style = "Normal"
